$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.033.13"
$ws.Range("E2").Value = "  +4.77%  "
$ws.Range("D3").Value = "2.242.17"
$ws.Range("E3").Value = "  +4.07%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("E5").Value = "  +4.06%  "
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "0.615"
$r.Style = "Normal"
$ws.Range("E6").Value = "  +2.39%  "
$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = "75.03"
$r.Style = "Normal"
$ws.Range("E7").Value = "  +9.39%  "
$ws.Range("E8").Value = "  -0.16%  "
$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = "0.604"
$r.Style = "Normal"
$ws.Range("E9").Value = "  +7.86%  "
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = "41.24"
$r.Style = "Normal"
$ws.Range("E10").Value = "  +8.26%  "
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "0.0928"
$r.Style = "Normal"
$ws.Range("E11").Value = "  +3.44%  "
$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = "6.90"
$r.Style = "Normal"
$ws.Range("E12").Value = "  +4.49%  "
$ws.Range("E13").Value = "  +1.99%  "
$ws.Range("D14").Value = "2.581.20"
$ws.Range("E14").Value = "  +3.99%  "
$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = "14.57"
$r.Style = "Normal"
$ws.Range("E15").Value = "  +2.08%  "
$ws.Range("D16").Value = "2.255.95"
$ws.Range("E16").Value = "  +5.84%  "
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = "0.791"
$r.Style = "Normal"
$ws.Range("E17").Value = "  +1.89%  "
$ws.Range("D18").Value = "42.951.93"
$ws.Range("E18").Value = "  +5.00%  "
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = "0.0000104"
$r.Style = "Normal"
$ws.Range("E19").Value = "  +6.19%  "
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "71.14"
$r.Style = "Normal"
$ws.Range("E20").Value = "  +2.42%  "
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = "5.98"
$r.Style = "Normal"
$ws.Range("E21").Value = "  +4.65%  "
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "9.80"
$r.Style = "Normal"
$ws.Range("E22").Value = "  +4.85%  "
$ws.Range("B23").Value = "ImmutableX"
$ws.Range("C23").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = "2.19"
$r.Style = "Normal"
$ws.Range("E23").Value = "  +18.75%  "
$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = "229.43"
$r.Style = "Normal"
$ws.Range("E24").Value = "  +2.86%  "
$ws.Range("E25").Value = "  +0.03%  "
$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = "10.84"
$r.Style = "Normal"
$ws.Range("E26").Value = "  +2.92%  "
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = "3.43"
$r.Style = "Normal"
$ws.Range("E27").Value = "  +7.84%  "
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = "2.23"
$r.Style = "Normal"
$ws.Range("E28").Value = "  +3.08%  "
$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = "38.57"
$r.Style = "Normal"
$ws.Range("E29").Value = "  +29.34%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = "171.66"
$r.Style = "Normal"
$ws.Range("E30").Value = "  +2.09%  "
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = "2.13"
$r.Style = "Normal"
$ws.Range("E31").Value = "  -1.23%  "
$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = "20.25"
$r.Style = "Normal"
$ws.Range("E32").Value = "  +3.45%  "
$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = "0.0799"
$r.Style = "Normal"
$ws.Range("E33").Value = "  +7.15%  "
$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = "5.29"
$r.Style = "Normal"
$ws.Range("E34").Value = "  +5.71%  "
$ws.Range("E35").Value = "  +2.13%  "
$ws.Range("E36").Value = "  +7.87%  "
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = "4.36"
$r.Style = "Normal"
$ws.Range("E37").Value = "  +7.98%  "
$ws.Range("E38").Value = "  +20.69%  "
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = "12.99"
$r.Style = "Normal"
$ws.Range("E39").Value = "  +14.98%  "
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "2.12"
$r.Style = "Normal"
$ws.Range("E40").Value = "  +4.27%  "
$ws.Range("E41").Value = "  +11.67%  "
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "5.43"
$r.Style = "Normal"
$ws.Range("E42").Value = "  +3.22%  "
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "104.84"
$r.Style = "Normal"
$ws.Range("E43").Value = "  +9.41%  "
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = "59.23"
$r.Style = "Normal"
$ws.Range("E44").Value = "  +4.14%  "
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = "8.66"
$r.Style = "Normal"
$ws.Range("E45").Value = "  +6.37%  "
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "0.477"
$r.Style = "Normal"
$ws.Range("E46").Value = "  +31.68%  "
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "0.0986"
$r.Style = "Normal"
$ws.Range("E47").Value = "  +3.75%  "
$ws.Range("E48").Value = "  +12.71%  "
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = "1.10"
$r.Style = "Normal"
$ws.Range("E49").Value = "  +3.72%  "
$ws.Range("E50").Value = "  +4.23%  "
$ws.Range("B51").Value = "HuobiToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = "2.68"
$r.Style = "Normal"
$ws.Range("E51").Value = "  +2.99%  "
